# Build the context from the variables used in the Excel file.
#
# - Sheet2 is populated with the truth-table combinations for the two
#   boolean variables (P, Q) that are used elsewhere in the workbook
#   (mirrors the P/Q truth table already present on Sheet1).
# - Sheet3 was an unused, empty scratch sheet and is removed.
# - Sheet2 becomes the active/selected sheet (with B4 selected), since
#   that's now where the context is built.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Sheet3 is unused (empty) - drop it, the context no longer needs it.
[void]$wb.Worksheets.Item("Sheet3").Delete()

# Build the context table on Sheet2 from the P/Q variables: every
# combination of the two booleans used in the truth table.
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Cells.Item(1, 1).Value = 0
$ws2.Cells.Item(1, 2).Value = 0

$ws2.Cells.Item(2, 1).Value = 0
$ws2.Cells.Item(2, 2).Value = 1

$ws2.Cells.Item(3, 1).Value = 1
$ws2.Cells.Item(3, 2).Value = 0

$ws2.Cells.Item(4, 1).Value = 1
$ws2.Cells.Item(4, 2).Value = 1

# Sheet1 rows keep their explicit (custom) row height now that the
# context sheet carries its own formatting.
$ws1 = $wb.Worksheets.Item("Sheet1")
for ($r = 1; $r -le 5; $r++) {
    $ws1.Rows.Item($r).RowHeight = 12.1
}

# Sheet2 is now the sheet the user is working from.
[void]$ws2.Activate()
[void]$ws2.Range("B4").Select()
